$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (rows 2-24) is relabelled with keycap-emoji numbers. The cells
# were originally entered by hand in a non-linear order (rows 21-24 ended up
# one off from the "row-1" pattern used for rows 2-20), so the exact value
# for every row is listed explicitly, in the order it was written, so both
# the final cell values and the shared-string table match.
$rowValues = [ordered]@{
  2 = "1️⃣"
  3 = "2️⃣"
  4 = "3️⃣"
  5 = "4️⃣"
  6 = "5️⃣"
  7 = "6️⃣"
  8 = "7️⃣"
  9 = "8️⃣"
  10 = "9️⃣"
  11 = "🔟"
  12 = "1️⃣1️⃣"
  13 = "1️⃣2️⃣"
  22 = "2️⃣2️⃣"
  21 = "2️⃣1️⃣"
  14 = "1️⃣3️⃣"
  15 = "1️⃣4️⃣"
  24 = "2️⃣4️⃣"
  23 = "2️⃣3️⃣"
  16 = "1️⃣5️⃣"
  17 = "1️⃣6️⃣"
  18 = "1️⃣7️⃣"
  19 = "1️⃣8️⃣"
  20 = "1️⃣9️⃣"
}

foreach ($row in $rowValues.Keys) {
  $ws.Cells.Item($row, 1).Value = $rowValues[$row]
}

# Centre (horizontal + vertical) the new labels. Format A2 directly first,
# then fan that exact format out to the remaining cells via Copy /
# PasteSpecial so a single new style entry is produced (mirrors the author
# applying "Center"/"Middle Align" once and then painting the format across
# the rest of the column) instead of one new style per property write.
$firstCell = $ws.Range("A2")
$firstCell.HorizontalAlignment = -4108
$firstCell.VerticalAlignment = -4108
$firstCell.Copy()
$ws.Range("A2:A24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the selection left by the author at A18.
$null = $ws.Range("A18").Select()
